# "logboek voor vandaag aangepast"
#
# Douwe's log sheet is extended with today's (Saturday's) work session:
#  - Friday's entry (row 10) gets its note text updated.
#  - A new time entry is logged in row 11 (14:30 - 16:00) describing the
#    "snakes and ladders" ascii-art puzzle work, with the note cell
#    word-wrapped.
#  - The day label "zaterdag" is added next to the new entry (row 12,
#    mirroring how earlier day labels sit one row below the first entry
#    of that day).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Douwe")

# --- Friday entry (row 10): update the note text -------------------------
# Update this cell first so the edited string reuses its existing
# shared-string slot instead of bumping every later index.
$ws.Range("D10").Value = "Functie gemaakt om eten te kunnen eten en begonnen met het maken van ruimte 9"

# --- Saturday entry (row 11): new start/end time + note ------------------
$ws.Range("A11").Value = 0.60416666666666663   # 14:30
$ws.Range("B11").Value = 0.66666666666666663   # 16:00

$ws.Range("D11").Value = "Op papier eerst een ''snakes and ladders'' ascii art gemaakt om het vervolgens in python uit te printen. `nHet is voor een klimpuzzel in ruimte 9`nOp deze website is het gedaan https://asciiflow.com/#/"
$ws.Range("D11").WrapText = $true
$ws.Rows("11").RowHeight = 45

# --- Day label for the new entry ------------------------------------------
$ws.Range("E12").Value = "zaterdag"

# --- Leave the cursor where the author left it ----------------------------
$ws.Range("E12").Select() | Out-Null
